$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row for the "Display touch Raspberry Pi 720 x 1280" item ---
# It goes right after "Buzzer 5V 12mm" (row 7) and before "Raspberry Pi 4GB Ram"
# (old row 8), so insert above the old row 8. This shifts rows 8-12 down to
# 9-13 and Excel auto-extends the SUM() range used by the Custo Total formula.
$ws.Rows.Item(8).Insert()

# The freshly inserted row has no formatting of its own yet; give it the same
# look as the row above (C7:D7 - "Buzzer 5V 12mm"), which is what Excel does
# when a new row is inserted in the middle of a formatted table.
$ws.Range("C7:D7").Copy()
$ws.Range("C8:D8").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C8:D8").Borders.Item(8).LineStyle = -4142  # xlEdgeTop -> none (matches body rows)
$excel.CutCopyMode = 0

# Fill in the new component.
$ws.Range("C8").Value = "Display touch Raspberry Pi 720 x 1280"
$ws.Range("D8").Value = 759.05

# De-duplicate the header-row style: before the edit, cellXfs held two
# identical entries for the "Componente"/"Custo" header cells; make C6:D6
# point at the same (first) style entry instead of the redundant one.
$ws.Range("C5:D5").Copy()
$ws.Range("C6:D6").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C6").Value = "Componente"
$ws.Range("D6").Value = "Custo"
$excel.CutCopyMode = 0

# Column D needs to widen slightly to fit the new, longer currency value.
$ws.Columns.Item(4).ColumnWidth = 11.33

$ws.Range("D14").Select()
